$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MetadataDict")

# --- Update vocabulary: rename MeasurementOrFact "termID"/"term" rows ---
$ws.Range("B33").Value = "measurmentTypeID"
$ws.Range("B35").Value = "measurementType"

# --- Filter the MeasurementOrFact table down to just the MeasurementOrFact rows ---
$rng = $ws.Range("A1:O56")
$rng.AutoFilter(1, @("MeasurementOrFact"))

# --- Column A best-fit width for the now-visible (narrower) data ---
$col = $ws.Columns.Item(1)
$col.ColumnWidth = 18.17

# --- Move the active selection to A31 ---
$ws.Activate()
$ws.Range("A31").Select()

Write-Output "done"
